$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Coin (B) and Link (C) columns are plain text and never look numeric,
# so they can be set directly. Price (D) and Volume (E) columns contain
# values that look numeric/date-like to Excel (e.g. "3.632.19", "0.490",
# "  +0.81%  ") so each is force-formatted as Text first to prevent COM
# from silently reinterpreting/rounding them as numbers.

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "66.754.77"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +0.81%  "

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.632.19"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +1.99%  "

# Row 4
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  +0.03%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "612.52"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +1.03%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "150.93"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  +4.03%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "3.631.52"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +2.01%  "

# Row 8
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  +0.04%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.490"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  -0.06%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.138"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  +1.01%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "8.02"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  +0.89%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.418"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  +0.99%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.251.07"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  +2.07%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.0000210"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  +1.61%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "30.23"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  +0.48%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.645.26"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  +2.42%  "

# Row 17
$ws.Range("B17").Value = "TRON"
$ws.Range("C17").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.118"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  +2.09%  "

# Row 18
$ws.Range("B18").Value = "WrappedBTC"
$ws.Range("C18").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "66.898.66"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  +0.90%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "11.70"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +2.30%  "

# Row 20
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  +2.75%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "15.16"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  +2.61%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "428.11"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  -0.24%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.618"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  +0.42%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "79.05"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  -0.05%  "

# Row 25
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  +0.06%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.0000124"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  +5.04%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "8.45"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  +6.88%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "9.61"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  +5.46%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.54"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  +1.21%  "

# Row 30
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  -0.05%  "

# Row 31
$ws.Range("B31").Value = "Kaspa"
$ws.Range("C31").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.162"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  +5.12%  "

# Row 32
$ws.Range("B32").Value = "Fetch.AI"
$ws.Range("C32").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.49"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  +2.49%  "

# Row 33
$ws.Range("B33").Value = "RenzoRestakedETH"
$ws.Range("C33").Value = "https://coinranking.com/coin/lKlJ_MC5M+renzorestakedeth-ezeth"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.630.55"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  +2.08%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "25.46"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  -0.29%  "

# Row 35
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  +0.34%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "5.69"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  +1.93%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.72"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  -1.10%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "178.00"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  +1.84%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0865"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  +2.20%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "5.25"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  +0.94%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.901"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  +0.50%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.91"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  -3.23%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "46.25"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  +0.07%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.61"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  +9.18%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.00"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  +0.04%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "25.18"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  -0.57%  "

# Row 48
$ws.Range("B48").Value = "EnergySwap"
$ws.Range("C48").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "24.13"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  +3.04%  "

# Row 49
$ws.Range("B49").Value = "ONDO"
$ws.Range("C49").Value = "https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.16"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  -3.49%  "

# Row 50
$ws.Range("B50").Value = "SuiNetwork"
$ws.Range("C50").Value = "https://coinranking.com/coin/3xJluUMvp+suinetwork-sui"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.981"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  +4.64%  "

# Row 51
$ws.Range("B51").Value = "Cosmos"
$ws.Range("C51").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "7.22"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  +0.99%  "
